$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite data rows 1-28 (header + first 27 data rows) ---
$ws.Range("A1").Value = 'file'
$ws.Range("B1").Value = 'varcol'
$ws.Range("C1").Value = 'datecol'
$ws.Range("D1").Value = 'rowstart'
$ws.Range("E1").Value = 'group'
$ws.Range("F1").Value = 'transformation'
$ws.Range("G1").Value = 'source'
$ws.Range("A2").Value = 'PCEPILFE.xls'
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = 'inflation'
$ws.Range("F2").Value = 'D12'
$ws.Range("G2").Value = 'FRED'
$ws.Range("A3").Value = 'GBRCPIALLMINMEI.xls'
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 'inflation'
$ws.Range("F3").Value = 'D12'
$ws.Range("G3").Value = 'FRED'
$ws.Range("A4").Value = 'JPNCPIALLMINMEI.xls'
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 10
$ws.Range("E4").Value = 'inflation'
$ws.Range("F4").Value = 'D12'
$ws.Range("G4").Value = 'FRED'
$ws.Range("A5").Value = 'EXUSUK.xls'
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 10
$ws.Range("E5").Value = 'ex'
$ws.Range("F5").ClearContents()
$ws.Range("G5").Value = 'FRED'
$ws.Range("A6").Value = 'EXUSUK.xls'
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 10
$ws.Range("E6").Value = 'ex'
$ws.Range("F6").Value = 'LN'
$ws.Range("G6").Value = 'FRED'
$ws.Range("A7").Value = 'EXJPUS.xls'
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 10
$ws.Range("E7").Value = 'ex'
$ws.Range("F7").ClearContents()
$ws.Range("G7").Value = 'FRED'
$ws.Range("A8").Value = 'EXJPUS.xls'
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 10
$ws.Range("E8").Value = 'ex'
$ws.Range("F8").Value = 'LN'
$ws.Range("G8").Value = 'FRED'
$ws.Range("A9").Value = 'DBAA.xls'
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 10
$ws.Range("E9").Value = 'rate'
$ws.Range("F9").ClearContents()
$ws.Range("G9").Value = 'FRED'
$ws.Range("A10").Value = 'GSPTSE.xlsx'
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 'stock'
$ws.Range("F10").ClearContents()
$ws.Range("G10").Value = 'Yahoo'
$ws.Range("A11").Value = 'GSPTSE.xlsx'
$ws.Range("B11").Value = 7
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 'stock'
$ws.Range("F11").Value = 'LND1'
$ws.Range("G11").Value = 'Yahoo'
$ws.Range("A12").Value = 'INDPRO.xls'
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 10
$ws.Range("E12").Value = 'output'
$ws.Range("F12").ClearContents()
$ws.Range("G12").Value = 'FRED'
$ws.Range("A13").Value = 'CPILFESL.xls'
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 10
$ws.Range("E13").Value = 'inflation'
$ws.Range("F13").ClearContents()
$ws.Range("G13").Value = 'FRED'
$ws.Range("A14").Value = 'CPILFESL.xls'
$ws.Range("B14").Value = 2
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 10
$ws.Range("E14").Value = 'inflation'
$ws.Range("F14").Value = 'D12'
$ws.Range("G14").Value = 'FRED'
$ws.Range("A15").Value = 'N225.xlsx'
$ws.Range("B15").Value = 5
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 'stock'
$ws.Range("F15").ClearContents()
$ws.Range("G15").Value = 'Yahoo'
$ws.Range("A16").Value = 'N225.xlsx'
$ws.Range("B16").Value = 7
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 'stock'
$ws.Range("F16").Value = 'LND1'
$ws.Range("G16").Value = 'Yahoo'
$ws.Range("A17").Value = 'SP500.xlsx'
$ws.Range("B17").Value = 5
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 'stock'
$ws.Range("F17").ClearContents()
$ws.Range("G17").Value = 'Yahoo'
$ws.Range("A18").Value = 'SP500.xlsx'
$ws.Range("B18").Value = 7
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 'stock'
$ws.Range("F18").Value = 'LND1'
$ws.Range("G18").Value = 'Yahoo'
$ws.Range("A19").Value = 'series-060822.xlsx'
$ws.Range("B19").Value = 3
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 'unemp'
$ws.Range("F19").ClearContents()
$ws.Range("G19").Value = 'UK ONS'
$ws.Range("A20").Value = 'series-060822.xlsx'
$ws.Range("B20").Value = 4
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = 'unemp'
$ws.Range("F20").Value = 'D12'
$ws.Range("G20").Value = 'UK ONS'
$ws.Range("A21").Value = 'LRUN64TTJPM156S.xls'
$ws.Range("B21").Value = 2
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = 'unemp'
$ws.Range("F21").ClearContents()
$ws.Range("G21").Value = 'FRED'
$ws.Range("A22").Value = 'LRUN64TTJPM156S.xls'
$ws.Range("B22").Value = 3
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 10
$ws.Range("E22").Value = 'unemp'
$ws.Range("F22").Value = 'D12'
$ws.Range("G22").Value = 'FRED'
$ws.Range("A23").Value = 'DEXCAUS.xls'
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 10
$ws.Range("E23").Value = 'ex'
$ws.Range("F23").ClearContents()
$ws.Range("G23").Value = 'FRED'
$ws.Range("A24").Value = 'DEXCAUS.xls'
$ws.Range("B24").Value = 2
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = 'ex'
$ws.Range("F24").Value = 'LN'
$ws.Range("G24").Value = 'FRED'
$ws.Range("A25").Value = 'UNRATE.xls'
$ws.Range("B25").Value = 2
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 'unemp'
$ws.Range("F25").ClearContents()
$ws.Range("G25").Value = 'FRED'
$ws.Range("A26").Value = 'UNRATE.xls'
$ws.Range("B26").Value = 3
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 'unemp'
$ws.Range("F26").Value = 'D12'
$ws.Range("G26").Value = 'FRED'
$ws.Range("A27").Value = 'LRUNTTTTCAM156S.xls'
$ws.Range("B27").Value = 2
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 10
$ws.Range("E27").Value = 'unemp'
$ws.Range("F27").ClearContents()
$ws.Range("G27").Value = 'FRED'
$ws.Range("A28").Value = 'LRUNTTTTCAM156S.xls'
$ws.Range("B28").Value = 3
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 10
$ws.Range("E28").Value = 'unemp'
$ws.Range("F28").Value = 'D12'
$ws.Range("G28").Value = 'FRED'

# --- Apply AutoFilter over A1:G28 before adding the trailing rows ---
$ws.Range("A1:G28").AutoFilter(1) | Out-Null
$fd = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$G`$28")
$fd.Visible = $false

# --- Add the two new trailing rows (FTSE100.xlsx) after the filter range is set ---
$ws.Range("A29").Value = 'FTSE100.xlsx'
$ws.Range("B29").Value = 11
$ws.Range("C29").Value = 9
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 'stock'
$ws.Range("F29").ClearContents()
$ws.Range("G29").Value = 'WSJ'
$ws.Range("A30").Value = 'FTSE100.xlsx'
$ws.Range("B30").Value = 12
$ws.Range("C30").Value = 9
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 'stock'
$ws.Range("F30").Value = 'LND1'
$ws.Range("G30").Value = 'WSJ'

# --- Freeze header row and set active selection/cell ---
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E4").Select() | Out-Null
